$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "Nombre Completo: Juan Alvarez" "Nombre Completo: Carlos Fontaner"
Replace-Text "Número Exterior: COLONIA" "Número Exterior: COLONIA2"
Replace-Text "Ciudad: PUEBLAYORK" "Ciudad: PUEBLA"
Replace-Text "Salario Diario: `$5000.00" "Salario Diario: `$10000.00"
Replace-Text "Fecha de Ingreso: 2024-11-04" "Fecha de Ingreso: 2024-11-03"
Replace-Text "Correo Electrónico: juanito@gmail.com" "Correo Electrónico: carlitos@gmail.com"
Replace-Text "comenzando el 2024-11-04" "comenzando el 2024-11-03"
Replace-Text "salario diario de `$5000.00" "salario diario de `$10000.00"
